$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F ("TÌNH TRẠNG" / status) gets marked "Hoàn thành" for several rows ---
$ws.Range("F19").Value = "Hoàn thành"
$ws.Range("F20").Value = "Hoàn thành"
$ws.Range("F23").Value = "Hoàn thành"
$ws.Range("F24").Value = "Hoàn thành"
$ws.Range("F25").Value = "Hoàn thành"

# --- Rows 27-29 ("Viết tài liệu..." section): fill start/end dates, member, status ---
$ws.Range("C27").Value = "18/10/23"
$ws.Range("D27").Value = "25/10/23"
$ws.Range("E27").Value = "Bùi Phi Long"
$ws.Range("F27").Value = "Hoàn thành"

$ws.Range("C28").Value = "18/10/23"
$ws.Range("D28").Value = "25/10/23"
$ws.Range("E28").Value = "Nguyễn Phạm Nhật Minh"
$ws.Range("F28").Value = "Hoàn thành"

$ws.Range("C29").Value = "18/10/23"
$ws.Range("D29").Value = "25/10/23"
$ws.Range("E29").Value = "Nguyễn Phạm Nhật Minh"
$ws.Range("F29").Value = "Hoàn thành"

# --- Update the view: scroll down a bit and move the selection ---
$ws.Range("C31").Select()
